$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-18 Monday", "2025-08-19 Tuesday"),
    @("662÷2=331, 0", "981÷3=327, 0"),
    @("203÷4=50, 3", "132÷6=22, 0"),
    @("660÷6=110, 0", "702÷5=140, 2"),
    @("231÷9=25, 6", "465÷9=51, 6"),
    @("158÷7=22, 4", "792÷4=198, 0"),
    @("664÷8=83, 0", "774÷7=110, 4"),
    @("658÷6=109, 4", "201÷4=50, 1"),
    @("606÷8=75, 6", "588÷4=147, 0"),
    @("727÷8=90, 7", "790÷7=112, 6"),
    @("534÷6=89, 0", "129÷4=32, 1"),
    @("671÷3=223, 2", "216÷4=54, 0"),
    @("554÷7=79, 1", "187÷8=23, 3"),
    @("396÷3=132, 0", "168÷3=56, 0"),
    @("432÷4=108, 0", "697÷6=116, 1"),
    @("911÷9=101, 2", "917÷9=101, 8"),
    @("518÷6=86, 2", "177÷6=29, 3"),
    @("631÷6=105, 1", "643÷6=107, 1"),
    @("195÷4=48, 3", "483÷8=60, 3"),
    @("884÷4=221, 0", "764÷2=382, 0"),
    @("678÷4=169, 2", "800÷5=160, 0"),
    @("395÷3=131, 2", "221÷7=31, 4"),
    @("970÷9=107, 7", "971÷8=121, 3"),
    @("834÷6=139, 0", "694÷5=138, 4"),
    @("793÷2=396, 1", "507÷7=72, 3"),
    @("845÷8=105, 5", "717÷3=239, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
